$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new data row (row 17), mirroring the formatting of the prior
# last row (row 16) so the date column keeps its existing style.
$ws.Range("A16:N16").Copy($ws.Range("A17:N17"))

$ws.Cells.Item(17, 1).Value = 42625.885821759257
$ws.Cells.Item(17, 2).Value = 35
$ws.Cells.Item(17, 3).Value = 0
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(17, 9).Value = 0
$ws.Cells.Item(17, 10).Value = 0
$ws.Cells.Item(17, 11).Value = 0
$ws.Cells.Item(17, 12).Value = 0
$ws.Cells.Item(17, 13).Value = 0
$ws.Cells.Item(17, 14).Value = "Random"
